# Insert a new weekly price record as row 176, pushing the existing
# rows 176-215 down to 177-216 (mirrors the commit "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 176..215 down to 177..216 by inserting a new blank row at 176.
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new weekly record.
# (Same market / region / product as its neighbours; only the date, volume,
# weighted price and unit of sale differ from what used to be row 176.)
$ws.Cells.Item(176, 1).Value2  = 3
$ws.Cells.Item(176, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(176, 3).Value2  = "Coquimbo"
$ws.Cells.Item(176, 4).Value2  = 44855
$ws.Cells.Item(176, 5).Value2  = 5
$ws.Cells.Item(176, 6).Value2  = 100112026
$ws.Cells.Item(176, 7).Value2  = "Haba"
$ws.Cells.Item(176, 8).Value2  = "Sin especificar"
$ws.Cells.Item(176, 9).Value2  = "Primera"
$ws.Cells.Item(176, 10).Value2 = 115
$ws.Cells.Item(176, 11).Value2 = 8000
$ws.Cells.Item(176, 12).Value2 = 8500
$ws.Cells.Item(176, 13).Value2 = 8261
$ws.Cells.Item(176, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(176, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(176, 16).Value2 = 330
$ws.Cells.Item(176, 17).Value2 = 25
$ws.Cells.Item(176, 18).Value2 = "Hortaliza"
